$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "onda" is entered first so it lands before "degrau" in the shared-strings
# table (matches the authoring order baked into the target workbook).
$ws.Range("A5").Value = "onda"
$ws.Range("A5:A6").Merge()
$ws.Range("A5:A6").HorizontalAlignment = -4108
$ws.Range("A5:A6").VerticalAlignment = -4108

$ws.Range("A4").Value = "degrau"
$ws.Range("A7").Value = "degrau"

# The note that used to live at B14 moves up to O6 (same text, new spot).
$ws.Range("O6").Value = "Ymax e Tmax, sujeitos a erro"
$ws.Range("B14").Clear()

# Column B grew a bit to fit the new labels.
$ws.Columns("B:B").ColumnWidth = 8.86

$ws.Range("G7").Select()
